$wb = $excel.ActiveWorkbook

# --- SWERVING: add a new row (def456, sedan, red, July 2, 2016, 3:26pm) ---
$wsSwerving = $wb.Worksheets.Item("SWERVING")
$wsSwerving.Range("A3").Value = "def456"
$wsSwerving.Range("B3").Value = "sedan"
$wsSwerving.Range("C3").Value = "red"
$wsSwerving.Range("D3").Value = "July 2, 2016"
$wsSwerving.Range("E3").Value = "3:26pm"

# --- BEATING THE RED LIGHT: add a new row (abc, government, green, August 2, 2015, 8:00am) ---
$wsRedLight = $wb.Worksheets.Item("BEATING THE RED LIGHT")
$wsRedLight.Range("A2").Value = "abc"
$wsRedLight.Range("B2").Value = "government"
$wsRedLight.Range("C2").Value = "green"
$wsRedLight.Range("D2").Value = "August 2, 2015"
$wsRedLight.Range("E2").Value = "8:00am"

# --- COLOR CODING: update row 2 values and remove row 3 ---
$wsColorCoding = $wb.Worksheets.Item("COLOR CODING")
$wsColorCoding.Range("A2").Value = "abcd"
$wsColorCoding.Range("B2").Value = "suv"
$wsColorCoding.Range("C2").Value = "orange"
$wsColorCoding.Rows.Item(3).Delete()

# --- Update selections on each affected sheet ---
$wsSpeeding = $wb.Worksheets.Item("SPEEDING")
$wsSpeeding.Activate()
$wsSpeeding.Range("A3:E3").Select()

$wsSwerving.Activate()
$wsSwerving.Range("E4").Select()

$wsRedLight.Activate()
$wsRedLight.Range("C3").Select()

# --- COLOR CODING becomes the active/selected tab ---
$wsColorCoding.Activate()
$wsColorCoding.Range("A2").Select()
